$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new trade row (row 4) mirroring the structure/formatting of row 3.
$ws.Cells.Item(3, 1).Copy()
$ws.Cells.Item(4, 1).PasteSpecial(-4122)  # xlPasteFormats
$ws.Cells.Item(4, 1).Value = 42633.676689814813

$ws.Cells.Item(4, 2).Value = $false

$ws.Cells.Item(4, 3).Value = 9956.5
$ws.Cells.Item(4, 4).Value = 10000
$ws.Cells.Item(4, 5).Value = 108.67
$ws.Cells.Item(4, 6).Value = 107.73

$ws.Cells.Item(3, 7).Copy()
$ws.Cells.Item(4, 7).PasteSpecial(-4122)  # xlPasteFormats
$ws.Cells.Item(4, 7).Value = $false

$ws.Cells.Item(4, 8).Value = -0.87

$ws.Cells.Item(4, 9).Value = $false

$excel.CutCopyMode = $false
